$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.114.88'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").Value = '1.880.87'
$ws.Range("E3").Value = '  -0.85%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '313.51'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = '0.5073'
$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("D8").Value = '0.3867'
$ws.Range("E8").Value = '  -1.34%  '

$ws.Range("D9").Value = '0.09023'
$ws.Range("E9").Value = '  -3.72%  '

$ws.Range("D10").Value = '1.126'
$ws.Range("E10").Value = '  -0.46%  '

$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").Value = '6.372'
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '20.82'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.876.13'
$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '7.263'
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = '1.002'
$ws.Range("E15").Value = '  +0.12%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001114'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '91.47'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06625'
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '18.25'
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.130'
$ws.Range("E21").Value = '  -1.57%  '

$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.140.49'
$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '11.43'
$ws.Range("E23").Value = '  +0.86%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.271'
$ws.Range("E24").Value = '  -1.88%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.098.46'
$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.544'
$ws.Range("E26").Value = '  -2.95%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '20.83'
$ws.Range("E27").Value = '  -0.26%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '156.82'
$ws.Range("E28").Value = '  -0.16%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '127.16'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.1061'
$ws.Range("E30").Value = '  -0.28%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.065'
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.631'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.603'
$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").Value = '9.599'
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.06635'
$ws.Range("E35").Value = '  +0.39%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02412'
$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2191'
$ws.Range("E37").Value = '  +0.83%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.283'
$ws.Range("E38").Value = '  -1.23%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '1.215'
$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6423'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '11.49'
$ws.Range("E41").Value = '  +0.49%  '

$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").Value = '4.921'
$ws.Range("E42").Value = '  -1.38%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").Value = '0.6056'
$ws.Range("E44").Value = '  +1.15%  '

$ws.Range("D45").Value = '13.20'
$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("D46").Value = '1.277'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("E47").Value = '  -1.23%  '

$ws.Range("D48").Value = '1.250'
$ws.Range("E48").Value = '  +6.29%  '

$ws.Range("D49").Value = '2.007'
$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").Value = '121.53'
$ws.Range("E50").Value = '  -1.19%  '

$ws.Range("D51").Value = '79.66'
$ws.Range("E51").Value = '  +2.45%  '
